$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Delete()

$text = @'
questions = [
    {
        "title": "How can you increase the memory limit of WordPress?",
        "ques_type": 2,
        "options": [
            "define('WP_MEMORY_ NOEXCEED', '256M') in wp-setting.php",
            "define('WP_MEMORY_LIMIT', '256M') in wp-setting.php",
            "define('WP_MEMORY_LIMIT', '256M') in wp-config.php",
            "define('WP_MEMORY_NOEXCEED', '256M') in wp-config.php"
        ],
        "score": "define('WP_MEMORY_LIMIT', '256M') in wp-config.php"
    },
    {
        "title": "Hiding the version number of your Wordpress installation is a good security practice. How do you do it?",
        "ques_type": 2,
        "options": [
            "remove_action('wp_remove_version', 'wp_generator')",
            "remove_action('wp_version', 'wp_generator')",
            "remove_action('wp_head', 'wp_generator')",
            "remove_action('wp_head', 'wp_remove_version')"
        ],
        "score": "remove_action('wp_head', 'wp_generator')"
    },
    {
        "title": "What piece of code returns the directory path of the plugin?",
        "ques_type": 2,
        "options": [
            "&lt?php plugin_basename($file) ?&gt\n",
            "&lt?php plug_dir_path( $file ) ?&gt",
            "&lt?php plugin_info(($file) ?&gt",
            "&lt?phpplugin_content_dir($file) ?&gt"
        ],
        "score": "&lt?php plug_dir_path( $file ) ?&gt"
    },
    {
        "title": "What are the WordPress functions related to comments?",
        "ques_type": 15,
        "options": [
            "wp_allow_comment",
            "wp_remove_comment",
            "wp_count_comment",
            "wp_block_comment",
            "wp_delete_comment",
            "wp_publish_comment"
        ],
        "score": [
            "wp_allow_comment",
            "wp_count_comment",
            "wp_delete_comment"
        ]
    }
]
'@

$ws.Range("A1").Value = $text
